# Auto-generated Word COM-interop script (PowerShell-style) for the LOT2069.docx restructuring.
# The edit rearranges several content blocks between paragraphs/runs while each paragraph keeps
# its own formatting (pPr/rPr) fixed in place. We therefore replace the TEXT of each affected run
# with its target text. A two-phase placeholder swap avoids any accidental cross-matches while
# the rearrangement is applied.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        $preview = $old.Substring(0, [Math]::Min(60, $old.Length))
        throw "Replace failed, text not found: $preview"
    }
}

# --- Phase 1: tag each affected block with a unique placeholder ---
$old0 = @'
Transmitir aos alunos os conceitos básicos relacionados diretamente a engenharia de sistemas biológicos capacitando-os ao entendimento dos princípios de engenharia envolvidos em operações em larga escala, em sistemas com organismos vivos, ecossistemas e processos biológicos.Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, permitindo aos alunos estudar tópicos avançados em Engenharia de biossistemas, em uma abordagemvariável e multidisciplinar em temas relevantes a Engenharia.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos
'@
Replace-Exact $old0 "@@SLOT0@@"

$old1 = @'
Transmit to students the basic concepts directly related to the engineering of biological systems, enabling them to understand the engineering principles involved in large-scale operations, in systems with living organisms, ecosystems and biological processes.Present students with an insight into the potential and strategic applications of modern biotechnology, allowing students to study advanced topics in Biosystems Engineering, in variable and multidisciplinar approach in topics relevant to Engineering.Improve reasoning and awaken students’ critical spirit and creativity.
'@
Replace-Exact $old1 "@@SLOT1@@"

$old2 = @'
1814052 - Silvio Silverio da Silva
'@
Replace-Exact $old2 "@@SLOT2@@"

$old3 = @'
Introdução. Aspectos de Engenharia em processos fermentativos e enzimáticos envolvendo sistemas biológicos. Processos biotecnológicos de importância industrial. Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental.2 Análise de critérios de ampliação de escala em processos envolvendo sistemas biológicos. Introdução às técnicas de separação/purificação de produtos biotecnológicos.
'@
Replace-Exact $old3 "@@SLOT3@@"

$old4 = @'
Introduction. Engineering aspects in fermentative and enzymatic processes involving biological systems. Biotechnological processes of industrial importance. Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation.Analysis of scale up criteria in processes involving biological systems. Introduction to separation/purification techniques for biotechnological products.
'@
Replace-Exact $old4 "@@SLOT4@@"

$old5 = @'
- Introdução: importância dos bioprocessos e biossistemas e aplicações industriais.- Aspectos de Engenharia aplicados em Processos fermentativos e enzimáticos: características, biorreatores, operações, controle, sensores utilizados, aspectos cinéticos e modelagem de biossistemas.- Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.- Fundamentos de engenharia de bioprocessos aplicados aos biossistemas utilizando organismos vivos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado.- Análise de critérios de variação de escala em processos envolvendo sistemas biológicos.- Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental, exemplos práticos e estudo de casos.- Introdução às técnicas de separação/purificação de produtos biotecnológicos.
'@
Replace-Exact $old5 "@@SLOT5@@"

$old6 = @'
Os alunos serão avaliados formalmente por apresentação de trabalhos/estudos de casos e seminários aplicados durante o curso
'@
Replace-Exact $old6 "@@SLOT6@@"

$old7 = @'
Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5.
'@
Replace-Exact $old7 "@@SLOT7@@"

$old8 = @'
A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0.
'@
Replace-Exact $old8 "@@SLOT8@@"

$old9 = @'
1. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 1ª ed. - Edgard Blucher, 20012. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 2ª ed. - Edgard Blucher, 2020.3. AQUARONE, E. et al. Biotecnologia Industrial, vol. 4 – Biotecnologia na Produção de Alimentos - Edgard Blucher, 2001. 4. CASTILHO, L.R.; AUGUSTO, E.F.P.; MORAES, A. Tecnologia de Cultivo de Células Animais - de Biofármacos à Terapia Gênica. Roca, 2008.5. PESSOA JR, Adalberto et al. Biotecnologia farmacêutica: Aspectos sobre aplicação industrial. Editora Blucher, 2021.6. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 1. - Bioprocessos. Elsevier, 2017.7. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 2 - Alimentos. Elsevier, 2017.
'@
Replace-Exact $old9 "@@SLOT9@@"

# --- Phase 2: replace each placeholder with its final target text ---
$new0 = @'
Introdução. Aspectos de Engenharia em processos fermentativos e enzimáticos envolvendo sistemas biológicos. Processos biotecnológicos de importância industrial. Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental.2 Análise de critérios de ampliação de escala em processos envolvendo sistemas biológicos. Introdução às técnicas de separação/purificação de produtos biotecnológicos.
'@
Replace-Exact "@@SLOT0@@" $new0

$new1 = @'
Introduction. Engineering aspects in fermentative and enzymatic processes involving biological systems. Biotechnological processes of industrial importance. Engineering fundamentals applied to biosystems of importance in agriculture, medicine, biotechnology, biopharmaceuticals, industrial bioprocessing and environmental conservation.Analysis of scale up criteria in processes involving biological systems. Introduction to separation/purification techniques for biotechnological products.
'@
Replace-Exact "@@SLOT1@@" $new1

$new2 = @'
Transmitir aos alunos os conceitos básicos relacionados diretamente a engenharia de sistemas biológicos capacitando-os ao entendimento dos princípios de engenharia envolvidos em operações em larga escala, em sistemas com organismos vivos, ecossistemas e processos biológicos.Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, permitindo aos alunos estudar tópicos avançados em Engenharia de biossistemas, em uma abordagemvariável e multidisciplinar em temas relevantes a Engenharia.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos
'@
Replace-Exact "@@SLOT2@@" $new2

$new3 = @'
- Introdução: importância dos bioprocessos e biossistemas e aplicações industriais.- Aspectos de Engenharia aplicados em Processos fermentativos e enzimáticos: características, biorreatores, operações, controle, sensores utilizados, aspectos cinéticos e modelagem de biossistemas.- Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.- Fundamentos de engenharia de bioprocessos aplicados aos biossistemas utilizando organismos vivos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado.- Análise de critérios de variação de escala em processos envolvendo sistemas biológicos.- Fundamentos de engenharia aplicados aos biossistemas de importância na agricultura, medicina, biotecnologia, biofármacos, bioprocessamento industrial e conservação ambiental, exemplos práticos e estudo de casos.- Introdução às técnicas de separação/purificação de produtos biotecnológicos.
'@
Replace-Exact "@@SLOT3@@" $new3

$new4 = @'
Transmit to students the basic concepts directly related to the engineering of biological systems, enabling them to understand the engineering principles involved in large-scale operations, in systems with living organisms, ecosystems and biological processes.Present students with an insight into the potential and strategic applications of modern biotechnology, allowing students to study advanced topics in Biosystems Engineering, in variable and multidisciplinar approach in topics relevant to Engineering.Improve reasoning and awaken students’ critical spirit and creativity.
'@
Replace-Exact "@@SLOT4@@" $new4

$new5 = @'
Os alunos serão avaliados formalmente por apresentação de trabalhos/estudos de casos e seminários aplicados durante o curso
'@
Replace-Exact "@@SLOT5@@" $new5

$new6 = @'
Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5.
'@
Replace-Exact "@@SLOT6@@" $new6

$new7 = @'
A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0.
'@
Replace-Exact "@@SLOT7@@" $new7

$new8 = @'
1. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 1ª ed. - Edgard Blucher, 20012. LIMA, U.A. et al. Biotecnologia Industrial, vol. 3 - Processos Fermentativos e Enzimáticos - 2ª ed. - Edgard Blucher, 2020.3. AQUARONE, E. et al. Biotecnologia Industrial, vol. 4 – Biotecnologia na Produção de Alimentos - Edgard Blucher, 2001. 4. CASTILHO, L.R.; AUGUSTO, E.F.P.; MORAES, A. Tecnologia de Cultivo de Células Animais - de Biofármacos à Terapia Gênica. Roca, 2008.5. PESSOA JR, Adalberto et al. Biotecnologia farmacêutica: Aspectos sobre aplicação industrial. Editora Blucher, 2021.6. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 1. - Bioprocessos. Elsevier, 2017.7. NASCIMENTO, R. et al. Microbiologia Industrial, vol. 2 - Alimentos. Elsevier, 2017.
'@
Replace-Exact "@@SLOT8@@" $new8

$new9 = @'
1814052 - Silvio Silverio da Silva
'@
Replace-Exact "@@SLOT9@@" $new9

Write-Output "done"
